# Auto-generated edit script applying scheduled runner updates to Famfrit_Profits sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 765
$ws.Range("I58").Value = 450.7143
$ws.Range("K58").Value = 1352.1429
$ws.Range("M58").Value = -1202.1429
$ws.Range("H138").Value = 6670804
$ws.Range("I138").Value = 1119.0555
$ws.Range("J138").Value = 10422501
$ws.Range("K138").Value = 3357.1665
$ws.Range("L138").Value = 31267503
$ws.Range("M138").Value = 1782.8335
$ws.Range("N138").Value = -31277783

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 181.66667
$ws.Range("I4").Value = 72.5
$ws.Range("K4").Value = 72.5
$ws.Range("M4").Value = 43.5
$ws.Range("H5").Value = 186.76471
$ws.Range("I5").Value = 186.76471
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 186.76471
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -74.76471000000001
$ws.Range("H61").Value = 20838982
$ws.Range("I61").Value = 33337476
$ws.Range("K61").Value = 33337476
$ws.Range("M61").Value = -33337264
$ws.Range("H110").Value = 12877.903
$ws.Range("I110").Value = 14412.5
$ws.Range("K110").Value = 14412.5
$ws.Range("M110").Value = -12367.5
$ws.Range("H122").Value = 3458.8696
$ws.Range("I122").Value = 2678.6
$ws.Range("J122").Value = 4059.077
$ws.Range("K122").Value = 8035.799999999999
$ws.Range("L122").Value = 12177.231
$ws.Range("M122").Value = -5585.799999999999
$ws.Range("N122").Value = -17077.231
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").ClearContents()
$ws.Range("N128").Value = 0
$ws.Range("H136").Value = 20838982
$ws.Range("I136").Value = 33337476
$ws.Range("K136").Value = 100012428
$ws.Range("M136").Value = -100009878

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 186.76471
$ws.Range("I4").Value = 186.76471
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 186.76471
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -71.76471000000001
$ws.Range("H22").Value = 465.5
$ws.Range("I22").Value = 465.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 465.5
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -292.5
$ws.Range("H94").Value = 1597.7693
$ws.Range("I94").Value = 244.8
$ws.Range("K94").Value = 244.8
$ws.Range("M94").Value = 206.2
$ws.Range("H134").Value = 2433.0356
$ws.Range("I134").Value = 2344.6428
$ws.Range("J134").Value = 2521.4285
$ws.Range("K134").Value = 7033.928400000001
$ws.Range("L134").Value = 7564.2855
$ws.Range("M134").Value = -4498.928400000001
$ws.Range("N134").Value = -12634.2855

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 9443.727999999999
$ws.Range("I22").Value = 9443.727999999999
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 9443.727999999999
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -9093.727999999999
$ws.Range("H31").Value = 3908.0588
$ws.Range("I31").Value = 2893.4119
$ws.Range("J31").Value = 4922.706
$ws.Range("K31").Value = 2893.4119
$ws.Range("L31").Value = 4922.706
$ws.Range("M31").Value = -2598.4119
$ws.Range("N31").Value = -5512.706
$ws.Range("H32").Value = 1255
$ws.Range("I32").Value = 1255
$ws.Range("K32").Value = 1255
$ws.Range("M32").Value = -939
$ws.Range("H34").Value = 3908.0588
$ws.Range("I34").Value = 2893.4119
$ws.Range("J34").Value = 4922.706
$ws.Range("K34").Value = 2893.4119
$ws.Range("L34").Value = 4922.706
$ws.Range("M34").Value = -2691.4119
$ws.Range("N34").Value = -5326.706
$ws.Range("H58").Value = 3737.4119
$ws.Range("I58").Value = 3701.4666
$ws.Range("K58").Value = 3701.4666
$ws.Range("M58").Value = -3498.4666
$ws.Range("H99").Value = 22813.133
$ws.Range("I99").Value = 25168.924
$ws.Range("J99").Value = 7500.5
$ws.Range("K99").Value = 25168.924
$ws.Range("L99").Value = 7500.5
$ws.Range("M99").Value = -23670.924
$ws.Range("N99").Value = -10496.5
$ws.Range("H126").Value = 22813.133
$ws.Range("I126").Value = 25168.924
$ws.Range("J126").Value = 7500.5
$ws.Range("K126").Value = 75506.772
$ws.Range("L126").Value = 22501.5
$ws.Range("M126").Value = -73036.772
$ws.Range("N126").Value = -27441.5
$ws.Range("H132").Value = 56957.758
$ws.Range("I132").Value = 82162.75999999999
$ws.Range("K132").Value = 246488.28
$ws.Range("M132").Value = -243958.28
$ws.Range("H134").Value = 2451.9333
$ws.Range("I134").Value = 2231.5833
$ws.Range("K134").Value = 6694.749899999999
$ws.Range("M134").Value = -4159.749899999999
$ws.Range("H136").Value = 3737.4119
$ws.Range("I136").Value = 3701.4666
$ws.Range("K136").Value = 11104.3998
$ws.Range("M136").Value = -8554.399800000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 127.36364
$ws.Range("I7").Value = 146.55556
$ws.Range("K7").Value = 439.66668
$ws.Range("M7").Value = -327.66668
$ws.Range("H46").Value = 848.75
$ws.Range("I46").Value = 965
$ws.Range("K46").Value = 2895
$ws.Range("M46").Value = -2804
$ws.Range("H131").Value = 1662.3334
$ws.Range("J131").Value = 1651.2572
$ws.Range("L131").Value = 4953.7716
$ws.Range("N131").Value = -15033.7716
$ws.Range("H133").Value = 3898.5
$ws.Range("I133").Value = 3898.5
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 11695.5
$ws.Range("L133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -6635.5
$ws.Range("H134").Value = 4399.9565
$ws.Range("I134").Value = 1455.5
$ws.Range("K134").Value = 4366.5
$ws.Range("M134").Value = 703.5
$ws.Range("H140").Value = 1217.3334
$ws.Range("I140").Value = 990.375
$ws.Range("K140").Value = 2971.125
$ws.Range("M140").Value = 2208.875

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2908.4211
$ws.Range("I97").Value = 2203.5
$ws.Range("J97").Value = 4882.2
$ws.Range("K97").Value = 2203.5
$ws.Range("L97").Value = 4882.2
$ws.Range("M97").Value = -1707.5
$ws.Range("N97").Value = -5874.2
$ws.Range("H122").Value = 3585.5833
$ws.Range("I122").Value = 3504.5
$ws.Range("J122").Value = 3666.6667
$ws.Range("K122").Value = 10513.5
$ws.Range("L122").Value = 11000.0001
$ws.Range("M122").Value = -8063.5
$ws.Range("N122").Value = -15900.0001
$ws.Range("H126").Value = 3645.3076
$ws.Range("I126").Value = 3083.1667
$ws.Range("J126").Value = 4127.143
$ws.Range("K126").Value = 9249.500100000001
$ws.Range("L126").Value = 12381.429
$ws.Range("M126").Value = -6779.500100000001
$ws.Range("N126").Value = -17321.429
$ws.Range("H132").Value = 1042.3636
$ws.Range("I132").Value = 1027.1111
$ws.Range("J132").Value = 1111
$ws.Range("K132").Value = 3081.3333
$ws.Range("L132").Value = 3333
$ws.Range("M132").Value = -551.3333000000002
$ws.Range("N132").Value = -8393

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1147.8889
$ws.Range("J22").Value = 1740.4
$ws.Range("L22").Value = 1740.4
$ws.Range("N22").Value = -2330.4
$ws.Range("H27").Value = 1147.8889
$ws.Range("J27").Value = 1740.4
$ws.Range("L27").Value = 1740.4
$ws.Range("N27").Value = -1954.4
$ws.Range("H56").Value = 22697.5
$ws.Range("I56").Value = 900
$ws.Range("J56").Value = 44495
$ws.Range("K56").Value = 900
$ws.Range("L56").Value = 44495
$ws.Range("M56").Value = -209
$ws.Range("N56").Value = -45877
$ws.Range("H132").Value = 10764.667
$ws.Range("I132").Value = 840.8182
$ws.Range("J132").Value = 26359.285
$ws.Range("K132").Value = 2522.4546
$ws.Range("L132").Value = 79077.855
$ws.Range("M132").Value = 7.545399999999972
$ws.Range("N132").Value = -84137.855
$ws.Range("H136").Value = 2003.8
$ws.Range("I136").Value = 1001.5357
$ws.Range("J136").Value = 4342.4165
$ws.Range("K136").Value = 3004.6071
$ws.Range("L136").Value = 13027.2495
$ws.Range("M136").Value = -454.6071000000002
$ws.Range("N136").Value = -18127.2495

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1900
$ws.Range("I96").Value = 1900
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1900
$ws.Range("L96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -527
$ws.Range("H100").Value = 71429784
$ws.Range("J100").Value = 550
$ws.Range("L100").Value = 1100
$ws.Range("N100").Value = -2182
$ws.Range("H107").Value = 336.125
$ws.Range("I107").Value = 336.125
$ws.Range("K107").Value = 1008.375
$ws.Range("M107").Value = 911.625
$ws.Range("H122").Value = 61430
$ws.Range("I122").Value = 73523.78999999999
$ws.Range("K122").Value = 220571.37
$ws.Range("M122").Value = -218121.37
$ws.Range("H132").Value = 2481.8276
$ws.Range("I132").Value = 2791
$ws.Range("K132").Value = 8373
$ws.Range("M132").Value = -5843
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").ClearContents()
$ws.Range("N133").Value = 0
$ws.Range("H136").Value = 7801.6665
$ws.Range("I136").Value = 2404
$ws.Range("J136").Value = 8476.375
$ws.Range("K136").Value = 7212
$ws.Range("L136").Value = 25429.125
$ws.Range("M136").Value = -4662
$ws.Range("N136").Value = -30529.125
$ws.Range("H137").Value = 133489
$ws.Range("J137").Value = 133489
$ws.Range("L137").Value = 133489
$ws.Range("N137").Value = -143689
$ws.Range("H140").Value = 107714
$ws.Range("J140").Value = 107714
$ws.Range("L140").Value = 107714
$ws.Range("N140").Value = -118074
